# Generate Report for Archive
# Swap the "c1b996f8-b300-4010-a124-dd7686a1d3e0" and "d8c59ed8-0736-4ece-91a1-8cc8bf810546"
# records (row 4 <-> row 5) on every sheet, and flip the d8c59ed8 record's status
# from "Ready for handoff" to "In Translation" to reflect its refreshed state.
#
# Only the cells whose value actually differs between row 4 and row 5 are
# rewritten (columns that already hold an identical/empty value in both rows
# are left untouched so existing blank cells aren't disturbed).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (columns A:G, hyperlinks live in column B) ----
$ws = $wb.Worksheets.Item("Overview")

$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$g4 = $ws.Range("G4").Value()
$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$g5 = $ws.Range("G5").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("G4").Value = $g5
$ws.Range("A5").Value = $a4
$ws.Range("B5").Value = $b4
$ws.Range("G5").Value = $g4

# Row 4 now holds the d8c59ed8 record -> its status moved forward to "In Translation"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\d8c59ed8-0736-4ece-91a1-8cc8bf810546.md"
    } elseif ($addr -eq '$B$5') {
        $hl.TextToDisplay = "e2e\c1b996f8-b300-4010-a124-dd7686a1d3e0.md"
    }
}

# ---- zh-cn / de-de sheets (columns A:P, hyperlinks live in column A) ----
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $a4 = $ws.Range("A4").Value()
    $g4 = $ws.Range("G4").Value()
    $h4 = $ws.Range("H4").Value()
    $a5 = $ws.Range("A5").Value()
    $g5 = $ws.Range("G5").Value()
    $h5 = $ws.Range("H5").Value()

    $ws.Range("A4").Value = $a5
    $ws.Range("G4").Value = $g5
    $ws.Range("H4").Value = $h5
    $ws.Range("A5").Value = $a4
    $ws.Range("G5").Value = $g4
    $ws.Range("H5").Value = $h4

    # Row 4 now holds the d8c59ed8 record -> its status moved forward to "In Translation"
    $ws.Range("C4").Value = "In Translation"

    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$4') {
            $hl.TextToDisplay = "d8c59ed8-0736-4ece-91a1-8cc8bf810546.md"
        } elseif ($addr -eq '$A$5') {
            $hl.TextToDisplay = "c1b996f8-b300-4010-a124-dd7686a1d3e0.md"
        }
    }
}
